$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5388969779014587
$ws.Range("B1").Value = 4.222954273223877
$ws.Range("C1").Value = 3.970008611679077
$ws.Range("D1").Value = 1.336737990379333
$ws.Range("E1").Value = 0.7902191877365112
